{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\n// 1) The author removed the stray words \"\u043d\u0430 \u043f\u043e\u043f\u0443\u043b\u044f\u0440\u0438\u0437\u0430\u0446\u0456\u0454\u044e \" from the\n// opening sentence:\n//   \"\u0406\u0437 \u0440\u043e\u0437\u0432\u0438\u0442\u043a\u043e\u043c \u043d\u0430 \u043f\u043e\u043f\u0443\u043b\u044f\u0440\u0438\u0437\u0430\u0446\u0456\u0454\u044e \u0411\u041f\u041b\u0410 \u0442\u0430 \" -> \"\u0406\u0437 \u0440\u043e\u0437\u0432\u0438\u0442\u043a\u043e\u043c \u0411\u041f\u041b\u0410 \u0442\u0430 \"\nconst search1 = context.document.body.search(\"\u0406\u0437 \u0440\u043e\u0437\u0432\u0438\u0442\u043a\u043e\u043c \u043d\u0430 \u043f\u043e\u043f\u0443\u043b\u044f\u0440\u0438\u0437\u0430\u0446\u0456\u0454\u044e \u0411\u041f\u041b\u0410 \u0442\u0430 \", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearch1.load(\"items\");\nawait context.sync();\n\nif (search1.items.length > 0) {\n  search1.items[0].insertText(\"\u0406\u0437 \u0440\u043e\u0437\u0432\u0438\u0442\u043a\u043e\u043c \u0411\u041f\u041b\u0410 \u0442\u0430 \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) The author removed the redundant \"\u0442\u0430 \u043a\u043e\u0448\u0442\u0438 \" from the sentence about\n// saving time on choosing a body/motor/etc.:\n//   \"\u0437\u0430\u043e\u0449\u0430\u0434\u0436\u0443\u0454 \u0447\u0430\u0441 \u0442\u0430 \u043a\u043e\u0448\u0442\u0438 \u043d\u0430 \u0432\u0438\u0431\u0456\u0440\" -> \"\u0437\u0430\u043e\u0449\u0430\u0434\u0436\u0443\u0454 \u0447\u0430\u0441 \u043d\u0430 \u0432\u0438\u0431\u0456\u0440\"\nconst search2 = context.document.body.search(\"\u0437\u0430\u043e\u0449\u0430\u0434\u0436\u0443\u0454 \u0447\u0430\u0441 \u0442\u0430 \u043a\u043e\u0448\u0442\u0438 \u043d\u0430 \u0432\u0438\u0431\u0456\u0440\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearch2.load(\"items\");\nawait context.sync();\n\nif (search2.items.length > 0) {\n  search2.items[0].insertText(\"\u0437\u0430\u043e\u0449\u0430\u0434\u0436\u0443\u0454 \u0447\u0430\u0441 \u043d\u0430 \u0432\u0438\u0431\u0456\u0440\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop edit script.\n# $word is the Application, $d (set below) is the ActiveDocument.\n\n$d = $word.ActiveDocument\n\n# 1) The author removed the stray words \"\u043d\u0430 \u043f\u043e\u043f\u0443\u043b\u044f\u0440\u0438\u0437\u0430\u0446\u0456\u0454\u044e \" from the\n# opening sentence:\n#   \"\u0406\u0437 \u0440\u043e\u0437\u0432\u0438\u0442\u043a\u043e\u043c \u043d\u0430 \u043f\u043e\u043f\u0443\u043b\u044f\u0440\u0438\u0437\u0430\u0446\u0456\u0454\u044e \u0411\u041f\u041b\u0410 \u0442\u0430 \" -> \"\u0406\u0437 \u0440\u043e\u0437\u0432\u0438\u0442\u043a\u043e\u043c \u0411\u041f\u041b\u0410 \u0442\u0430 \"\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"\u0406\u0437 \u0440\u043e\u0437\u0432\u0438\u0442\u043a\u043e\u043c \u043d\u0430 \u043f\u043e\u043f\u0443\u043b\u044f\u0440\u0438\u0437\u0430\u0446\u0456\u0454\u044e \u0411\u041f\u041b\u0410 \u0442\u0430 \"\n$find1.Replacement.Text = \"\u0406\u0437 \u0440\u043e\u0437\u0432\u0438\u0442\u043a\u043e\u043c \u0411\u041f\u041b\u0410 \u0442\u0430 \"\n$find1.Forward = $true\n$find1.Wrap = 1  # wdFindContinue\n$find1.Execute([ref]$find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)  # wdReplaceAll\n\n# 2) The author removed the redundant \"\u0442\u0430 \u043a\u043e\u0448\u0442\u0438 \" from the sentence about\n# saving time on choosing a body/motor/etc.:\n#   \"\u0437\u0430\u043e\u0449\u0430\u0434\u0436\u0443\u0454 \u0447\u0430\u0441 \u0442\u0430 \u043a\u043e\u0448\u0442\u0438 \u043d\u0430 \u0432\u0438\u0431\u0456\u0440\" -> \"\u0437\u0430\u043e\u0449\u0430\u0434\u0436\u0443\u0454 \u0447\u0430\u0441 \u043d\u0430 \u0432\u0438\u0431\u0456\u0440\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"\u0437\u0430\u043e\u0449\u0430\u0434\u0436\u0443\u0454 \u0447\u0430\u0441 \u0442\u0430 \u043a\u043e\u0448\u0442\u0438 \u043d\u0430 \u0432\u0438\u0431\u0456\u0440\"\n$find2.Replacement.Text = \"\u0437\u0430\u043e\u0449\u0430\u0434\u0436\u0443\u0454 \u0447\u0430\u0441 \u043d\u0430 \u0432\u0438\u0431\u0456\u0440\"\n$find2.Forward = $true\n$find2.Wrap = 1  # wdFindContinue\n$find2.Execute([ref]$find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)  # wdReplaceAll\n"}
